# Correction bug diffictulté + taverne temporium + affichage score ecran de victoire
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correction bug difficulté (colonne D = "Difficulté") ---
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D8").Value = 2

# --- Taverne temporium : correction de la valeur de test (colonne V, ligne 5) ---
$ws.Range("V5").Value = 40

# --- Affichage score ecran de victoire : la sélection se retrouve sur la
#     dernière colonne du score (Z5) ---
$ws.Range("Z5").Select()
